$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")

# Update "want-to-go" counts (F) for unaffected rows
$ws1.Cells.Item(4, 6).Value = 110
$ws1.Cells.Item(5, 6).Value = 1722
$ws1.Cells.Item(6, 6).Value = 3295
$ws1.Cells.Item(7, 6).Value = 965
$ws1.Cells.Item(8, 6).Value = 2134
$ws1.Cells.Item(9, 6).Value = 2057
$ws1.Cells.Item(10, 6).Value = 1070
$ws1.Cells.Item(11, 6).Value = 570
$ws1.Cells.Item(14, 6).Value = 361
$ws1.Cells.Item(16, 6).Value = 25
$ws1.Cells.Item(18, 6).Value = 144
$ws1.Cells.Item(19, 6).Value = 1509
$ws1.Cells.Item(21, 6).Value = 671
$ws1.Cells.Item(22, 6).Value = 559
$ws1.Cells.Item(23, 6).Value = 11986
$ws1.Cells.Item(24, 6).Value = 12003
$ws1.Cells.Item(26, 6).Value = 676

# Update venue text
$ws1.Cells.Item(11, 4).Value = "逸景路462号珠江国际纺织城d区6层 珠江时尚馆"

# Insert new row for the added event at row 27
$ws1.Rows.Item(27).Insert()
$ws1.Range("A26").Copy()
$ws1.Range("A27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Cells.Item(27, 1).Value = 26
$ws1.Cells.Item(27, 2).Value = "2024-06-09"
$ws1.Cells.Item(27, 3).Value = "广州·猎魔盛宴专场票·狂魔哥见面会专场票·珠三角COMIC WORLD次元世界动漫游戏嘉年华"
$ws1.Cells.Item(27, 4).Value = "南洲路139号 小洲云文化艺术创意园"
$ws1.Cells.Item(27, 5).Value = "2024.06.09 10:00-06.09 17:00"
$ws1.Cells.Item(27, 6).Value = 2
$ws1.Cells.Item(27, 7).Value = 238
$ws1.Cells.Item(27, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85441"
$ws1.Cells.Item(27, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/qDibxb9x1715096390466.jpeg"

# Update F/G for rows shifted down by the insertion
$ws1.Cells.Item(28, 6).Value = 3
$ws1.Cells.Item(29, 6).Value = 284
$ws1.Cells.Item(30, 6).Value = 1878
$ws1.Cells.Item(32, 6).Value = 504
$ws1.Cells.Item(28, 7).Value = 55
$ws1.Cells.Item(29, 7).Value = 54
$ws1.Cells.Item(30, 7).Value = 68

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")

# Update "want-to-go" counts (F) for unaffected rows
$ws4.Cells.Item(3, 6).Value = 68
$ws4.Cells.Item(6, 6).Value = 110
$ws4.Cells.Item(7, 6).Value = 1722
$ws4.Cells.Item(8, 6).Value = 3295
$ws4.Cells.Item(9, 6).Value = 965
$ws4.Cells.Item(10, 6).Value = 2134
$ws4.Cells.Item(11, 6).Value = 2057
$ws4.Cells.Item(12, 6).Value = 1070
$ws4.Cells.Item(13, 6).Value = 570
$ws4.Cells.Item(16, 6).Value = 361
$ws4.Cells.Item(18, 6).Value = 25
$ws4.Cells.Item(22, 6).Value = 144
$ws4.Cells.Item(23, 6).Value = 1509
$ws4.Cells.Item(25, 6).Value = 671
$ws4.Cells.Item(26, 6).Value = 559
$ws4.Cells.Item(27, 6).Value = 11986
$ws4.Cells.Item(28, 6).Value = 12003
$ws4.Cells.Item(30, 6).Value = 676

# Update venue text
$ws4.Cells.Item(13, 4).Value = "逸景路462号珠江国际纺织城d区6层 珠江时尚馆"

# Insert new row for the added event at row 31
$ws4.Rows.Item(31).Insert()
$ws4.Range("A30").Copy()
$ws4.Range("A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws4.Cells.Item(31, 1).Value = 30
$ws4.Cells.Item(31, 2).Value = "2024-06-09"
$ws4.Cells.Item(31, 3).Value = "广州·猎魔盛宴专场票·狂魔哥见面会专场票·珠三角COMIC WORLD次元世界动漫游戏嘉年华"
$ws4.Cells.Item(31, 4).Value = "南洲路139号 小洲云文化艺术创意园"
$ws4.Cells.Item(31, 5).Value = "2024.06.09 10:00-06.09 17:00"
$ws4.Cells.Item(31, 6).Value = 2
$ws4.Cells.Item(31, 7).Value = 238
$ws4.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85441"
$ws4.Cells.Item(31, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/qDibxb9x1715096390466.jpeg"

# Update F/G for rows shifted down by the insertion
$ws4.Cells.Item(32, 6).Value = 3
$ws4.Cells.Item(33, 6).Value = 284
$ws4.Cells.Item(34, 6).Value = 1878
$ws4.Cells.Item(38, 6).Value = 504
$ws4.Cells.Item(39, 6).Value = 9
$ws4.Cells.Item(32, 7).Value = 55
$ws4.Cells.Item(33, 7).Value = 54
$ws4.Cells.Item(34, 7).Value = 68

# ---- Sheet "演出" ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(7, 6).Value = 9

# ---- Sheet "本地生活" ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 68

